$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add B2
$ws.Range("B2").Value = "fsdfsd"

# Row 3: A3 removed ("fd12321"), B3 added
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "kjk"

# Row 4: B4 added
$ws.Range("B4").Value = "realy&"

# Row 5: B5 added
$ws.Range("B5").Value = "lol"

# Row 6: A6 removed ("/add_vendor"), B6 added
$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "gfdgfdg"

# Row 7: B7 added
$ws.Range("B7").Value = "добавь"

# Row 8: B8 added
$ws.Range("B8").Value = "fdsf"

# Row 9: new row
$ws.Range("A9").Value = "/add_vendor"
$ws.Range("B9").Value = "u"

# Row 10: new row
$ws.Range("A10").Value = "/addven"
$ws.Range("B10").Value = "f"

# Row 11: new row
$ws.Range("A11").Value = "/addven"
$ws.Range("B11").Value = "fd"

# Row 12: new row
$ws.Range("A12").Value = "/addven"

# Row 13: new row
$ws.Range("A13").Value = "/addven"

# Row 14: new row
$ws.Range("A14").Value = "/addven"

# Row 15: new row
$ws.Range("A15").Value = "/addven"

# Row 16: new row
$ws.Range("A16").Value = "/addven"

# Row 17: new row
$ws.Range("A17").Value = "gfdgfdgfd"

# Row 18: new row
$ws.Range("A18").Value = "жив?"

# Row 19: new row
$ws.Range("A19").Value = "lol"
